# Replace Group2 todo list with latest
# Updates the "Mret Hein" sheet's Revised Due By / Original Due By2 / Start By
# dates (columns D, E, F for rows 8-15) to the latest values, and moves the
# active selection to reflect where the author left off.
#
# Dates are written via Value2 (raw date serials) rather than Value/date
# strings so the existing custom date number-format on each cell (style
# indices already present in the workbook) is preserved instead of Excel
# re-guessing a generic date format for the new input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mret Hein")
$ws.Activate()

# Row 8 - Create Queries          (12/10/2025, 12/13/2025, 12/13/2025)
$ws.Range("D8").Value2 = 46001
$ws.Range("E8").Value2 = 46004
$ws.Range("F8").Value2 = 46004

# Row 9 - Create Markdown         (12/11/2025, 12/13/2025, 12/13/2025)
$ws.Range("D9").Value2 = 46002
$ws.Range("E9").Value2 = 46004
$ws.Range("F9").Value2 = 46004

# Row 10 - Query Creation         (12/14/2025, 12/14/2025, 12/14/2025)
$ws.Range("D10").Value2 = 46005
$ws.Range("E10").Value2 = 46005
$ws.Range("F10").Value2 = 46005

# Row 11 - Markdown Creation      (12/14/2025, 12/14/2025, 12/14/2025)
$ws.Range("D11").Value2 = 46005
$ws.Range("E11").Value2 = 46005
$ws.Range("F11").Value2 = 46005

# Row 12 - Error Assistance       (12/15/2025, 12/15/2025, 12/15/2025)
$ws.Range("D12").Value2 = 46006
$ws.Range("E12").Value2 = 46006
$ws.Range("F12").Value2 = 46006

# Row 13 - Format Review          (12/15/2025, 12/15/2025, 12/15/2025)
$ws.Range("D13").Value2 = 46006
$ws.Range("E13").Value2 = 46006
$ws.Range("F13").Value2 = 46006

# Row 14 - FeedBack Session       (12/16/2025, 12/16/2025, 12/16/2025)
$ws.Range("D14").Value2 = 46007
$ws.Range("E14").Value2 = 46007
$ws.Range("F14").Value2 = 46007

# Row 15 - Pair Recording         (12/16/2025, 12/16/2025, 12/16/2025)
$ws.Range("D15").Value2 = 46007
$ws.Range("E15").Value2 = 46007
$ws.Range("F15").Value2 = 46007

# Reflect the author's final cursor position on this sheet
$ws.Range("H18").Select()
